# Weekly update: a new price observation is inserted at row 86 (Fecha 44484),
# pushing the previously-existing rows 86..163 down by one row (87..164).
# Row 164 is a brand-new row that now holds the data that used to live in row 163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift rows 163 downto 87 (shift block), each row takes the current
#        (not-yet-modified) values of the row immediately above it.
#        Columns D (date), J, K, L, M, P are the only ones that vary row to row;
#        everything else (A,B,C,E,F,G,H,I,N,O,Q,R) is constant across this block,
#        so row 164 is seeded from row 163's static text/number columns below.

for ($r = 164; $r -ge 87; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value()   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($src, 10).Value() # J - Volumen
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($src, 11).Value() # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($src, 12).Value() # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($src, 13).Value() # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($src, 16).Value() # P - Precio $/Kg
}

# --- 2. Row 164 needs the constant columns too (it did not exist before).
$ws.Cells.Item(164, 1).Value = $ws.Cells.Item(163, 1).Value()   # A - Mercado ID
$ws.Cells.Item(164, 2).Value = $ws.Cells.Item(163, 2).Value()   # B - Mercado
$ws.Cells.Item(164, 3).Value = $ws.Cells.Item(163, 3).Value()   # C - Region
$ws.Cells.Item(164, 5).Value = $ws.Cells.Item(163, 5).Value()   # E - Codreg
$ws.Cells.Item(164, 6).Value = $ws.Cells.Item(163, 6).Value()   # F - Categoria ID
$ws.Cells.Item(164, 7).Value = $ws.Cells.Item(163, 7).Value()   # G - Categoria
$ws.Cells.Item(164, 8).Value = $ws.Cells.Item(163, 8).Value()   # H - Variedad
$ws.Cells.Item(164, 9).Value = $ws.Cells.Item(163, 9).Value()   # I - Calidad
$ws.Cells.Item(164, 14).Value = $ws.Cells.Item(163, 14).Value() # N - Unidad de comercializacion
$ws.Cells.Item(164, 15).Value = $ws.Cells.Item(163, 15).Value() # O - Origen
$ws.Cells.Item(164, 17).Value = $ws.Cells.Item(163, 17).Value() # Q - Kg o Unidades
$ws.Cells.Item(164, 18).Value = $ws.Cells.Item(163, 18).Value() # R - Clasificacion

# --- 3. Row 86 is the new observation - it keeps its existing static columns
#        but gets brand-new measured values.
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 4).Value = 44484
$ws.Cells.Item(86, 10).Value = 3300
$ws.Cells.Item(86, 11).Value = 400
$ws.Cells.Item(86, 12).Value = 500
$ws.Cells.Item(86, 13).Value = 450
$ws.Cells.Item(86, 16).Value = 900
